$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original header layout: A1=name, B1=phone, C1=email, D1=address
# New header layout:      A1=code, B1=name, C1=npwp, D1=email, E1=address, F1=phone
# Plus a new data row 2 with sample company values.

$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "npwp"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "address"
$ws.Range("F1").Value = "phone"

# Force text storage so numeric-looking strings (npwp/phone) keep their
# leading zeros instead of being coerced into numbers.
$ws.Range("A2:F2").NumberFormat = "@"

$ws.Range("A2").Value = "C001"
$ws.Range("B2").Value = "Company First"
$ws.Range("C2").Value = "01923821093"
$ws.Range("D2").Value = "company@test.com"
$ws.Range("E2").Value = "Jl. Soekarno Hatta"
$ws.Range("F2").Value = "08771939021"
